# ajuste en la carga de datos
# Adds a new delivery record in row 2 of the registro_entrega sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these numeric/date-looking values to be stored
# as literal text (matching how the source data was loaded), not
# auto-converted to a date serial / number.
$ws.Range("A2").Value = "'2025-08-01"
$ws.Range("B2").Value = "'16:36:47"

# fecha_devolucion / hora_devolucion are blank for this record. Touch a
# harmless no-op property so the cell is still materialized in the row
# (present, empty) without changing its appearance.
$ws.Range("C2").Font.Bold = $false
$ws.Range("D2").Font.Bold = $false

$ws.Range("E2").Value = "VIERNES"
$ws.Range("F2").Value = "INVESTIGACIÓN CONTABLE"
$ws.Range("G2").Value = "M-401"
$ws.Range("H2").Value = "MICHAEL  GONZALEZ PULGARIN"

$ws.Range("I2").Value = "'1036941224"

$ws.Range("J2").Value = "18:00 a 19:00"
$ws.Range("K2").Value = "Entregada"

# observaciones is blank for this record as well.
$ws.Range("L2").Font.Bold = $false
